$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 91
$ws.Cells.Item($row, 1).Value = "AI Engineer with GoLang"
$ws.Cells.Item($row, 2).Value = "https://www.dice.com/job-detail/e42b3e5d-fbf2-4c0b-9a91-367ff65deb1a"
$ws.Cells.Item($row, 3).Value = "Austin, Texas"
$ws.Cells.Item($row, 4).Value = "Third Party"
$ws.Cells.Item($row, 5).Value = "Depends on Experience"
$ws.Cells.Item($row, 6).Value = "NasTech Global, Inc."
